# Update cryptos list cell values per the Oct 28 2024 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.676.86"
$ws.Range("E2").Value = "'  +2.26%  "
$ws.Range("D3").Value = "'2.526.57"
$ws.Range("E3").Value = "'  +2.24%  "
$ws.Range("D5").Value = "'594.27"
$ws.Range("E5").Value = "'  +2.02%  "
$ws.Range("D6").Value = "'177.10"
$ws.Range("E6").Value = "'  +1.52%  "
$ws.Range("E7").Value = "'  +0.03%  "
$ws.Range("D8").Value = "'0.520"
$ws.Range("E8").Value = "'  +1.54%  "
$ws.Range("D9").Value = "'2.526.07"
$ws.Range("E9").Value = "'  +2.25%  "
$ws.Range("D10").Value = "'0.146"
$ws.Range("E10").Value = "'  +5.92%  "
$ws.Range("E11").Value = "'  -1.08%  "
$ws.Range("E12").Value = "'  +0.93%  "
$ws.Range("D13").Value = "'0.339"
$ws.Range("E13").Value = "'  +1.58%  "
$ws.Range("D14").Value = "'2.988.71"
$ws.Range("E14").Value = "'  +2.58%  "
$ws.Range("D15").Value = "'26.22"
$ws.Range("E15").Value = "'  +3.31%  "
$ws.Range("D16").Value = "'68.549.24"
$ws.Range("E16").Value = "'  +2.14%  "
$ws.Range("E17").Value = "'  +0.78%  "
$ws.Range("D18").Value = "'2.522.42"
$ws.Range("E18").Value = "'  +1.66%  "
$ws.Range("D19").Value = "'11.10"
$ws.Range("E19").Value = "'  +1.56%  "
$ws.Range("D20").Value = "'7.51"
$ws.Range("E20").Value = "'  +0.64%  "
$ws.Range("D21").Value = "'352.48"
$ws.Range("E21").Value = "'  +1.13%  "
$ws.Range("E22").Value = "'  +5.66%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = "'  -0.15%  "
$ws.Range("D24").Value = "'70.86"
$ws.Range("E24").Value = "'  +2.02%  "
$ws.Range("D25").Value = "'4.23"
$ws.Range("E25").Value = "'  +1.01%  "
$ws.Range("E26").Value = "'  -5.41%  "
$ws.Range("D27").Value = "'8.99"
$ws.Range("E27").Value = "'  -2.35%  "
$ws.Range("D28").Value = "'2.690.82"
$ws.Range("E28").Value = "'  +3.65%  "
$ws.Range("D29").Value = "'0.996"
$ws.Range("E29").Value = "'  -0.27%  "
$ws.Range("D30").Value = "'0.0₃0893"
$ws.Range("E30").Value = "'  -0.66%  "
$ws.Range("D31").Value = "'508.24"
$ws.Range("E31").Value = "'  +1.80%  "
$ws.Range("E32").Value = "'  +0.96%  "
$ws.Range("E33").Value = "'  +1.69%  "
$ws.Range("E34").Value = "'  +1.34%  "
$ws.Range("E35").Value = "'  +0.02%  "
$ws.Range("E36").Value = "'  +0.05%  "
$ws.Range("D37").Value = "'162.92"
$ws.Range("E37").Value = "'  +0.65%  "
$ws.Range("B38").Value = "'WhiteBITCoin"
$ws.Range("C38").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "'18.68"
$ws.Range("E38").Value = "'  +0.01%  "
$ws.Range("B39").Value = "'EthereumClassic"
$ws.Range("C39").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D39").Value = "'18.40"
$ws.Range("E39").Value = "'  +1.27%  "
$ws.Range("B40").Value = "'ImmutableX"
$ws.Range("C40").Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D40").Value = "'1.32"
$ws.Range("E40").Value = "'  -0.20%  "
$ws.Range("B41").Value = "'Stacks"
$ws.Range("C41").Value = "'https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.77"
$ws.Range("E41").Value = "'  +5.09%  "
$ws.Range("E42").Value = "'  -0.06%  "
$ws.Range("E43").Value = "'  +0.90%  "
$ws.Range("E44").Value = "'  -0.36%  "
$ws.Range("E45").Value = "'  +1.25%  "
$ws.Range("D46").Value = "'152.74"
$ws.Range("E46").Value = "'  +6.98%  "
$ws.Range("E47").Value = "'  +2.69%  "
$ws.Range("D48").Value = "'0.521"
$ws.Range("E48").Value = "'  +2.38%  "
$ws.Range("D49").Value = "'0.0₆0259"
$ws.Range("E49").Value = "'  +1.23%  "
$ws.Range("E50").Value = "'  +2.38%  "
$ws.Range("E51").Value = "'  -0.20%  "
